{"js": "// DTSCCI-248: use external short name rather than court venue name\n// Change the merge-field placeholder text in the document title block from\n// \"hearingLocation.site_name\" to \"hearingLocation.external_short_name\".\n// (Only the title-block occurrence changes; the other uses of\n// \"hearingLocation.site_name\" elsewhere in the document are left as-is.)\n\nconst body = context.document.body;\n\n// This exact phrase only occurs once in the whole document (in the bold,\n// centred title paragraph that also contains \"writtenByJudge\" and the\n// \"<<else>> Online Civil Claims\" fallback), so it unambiguously identifies\n// the paragraph to edit even though \"hearingLocation.site_name\" by itself\n// appears several more times later in the document (inside \"!=null\"\n// conditions) and must stay untouched.\nconst results = body.search(\n  \"<<cs_{writtenByJudge}>><<hearingLocation.site_name>><<else>> Online Civil Claims<<es_>>\",\n  { matchCase: true }\n);\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the title placeholder text, found \" +\n      results.items.length\n  );\n}\n\n// Replace the whole paragraph's text (grabbed via its range) so that Word\n// collapses the run back to a single clean run instead of leaving stray\n// proofing-error markers behind.\nconst paragraph = results.items[0].paragraphs.getFirst();\nconst paragraphRange = paragraph.getRange();\nparagraphRange.insertText(\n  \"<<cs_{writtenByJudge}>><<hearingLocation.external_short_name>><<else>> Online Civil Claims<<es_>>\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# DTSCCI-248: use external short name rather than court venue name\n# Change the merge-field placeholder text in the document title block from\n# \"hearingLocation.site_name\" to \"hearingLocation.external_short_name\".\n# (Only the title-block occurrence changes; the other uses of\n# \"hearingLocation.site_name\" elsewhere in the document - inside the\n# \"!=null\" conditional checks - are left as-is.)\n\n$d = $word.ActiveDocument\n\n# This exact phrase only occurs once in the whole document (in the bold,\n# centred title paragraph that also contains \"writtenByJudge\" and the\n# \"<<else>> Online Civil Claims\" fallback), so it unambiguously identifies\n# the text to edit even though \"hearingLocation.site_name\" by itself\n# appears several more times later in the document and must stay untouched.\n$searchText = \"<<cs_{writtenByJudge}>><<hearingLocation.site_name>><<else>> Online Civil Claims<<es_>>\"\n$replaceText = \"<<cs_{writtenByJudge}>><<hearingLocation.external_short_name>><<else>> Online Civil Claims<<es_>>\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, $true, $false, $replaceText, 1)\n\nif (-not $found) {\n    throw \"Could not find the title placeholder text to replace\"\n}\n"}
